# Weekly update: insert a new price record as the most recent entry (row 123)
# for "Feria Lagunitas de Puerto Montt - Apio", pushing the existing history
# down by one row (old row 123 -> 124, ... old row 156 -> 157).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 123; Excel shifts rows 123:156 down to 124:157
# and the sheet dimension grows to A1:R157 automatically.
$ws.Rows.Item(123).Insert()

# Seed the new row with the same layout as the row right below it (which was
# the former row 123), so constant columns (A,B,C,E,F,G,H,N,O,Q,R) and number
# formatting/styles (e.g. the date style on column D) come along for free.
$ws.Range("A124:R124").Copy()
$ws.Range("A123").PasteSpecial()

# Now overwrite the fields that differ for this new weekly record.
$ws.Range("D123").Value = 44508
$ws.Range("J123").Value = 25
$ws.Range("K123").Value = 10000
$ws.Range("L123").Value = 10000
$ws.Range("M123").Value = 10000
$ws.Range("P123").Value = 1667
